$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $range = $Sheet.Range($CellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = $origStyle
}

Set-TextValue $ws 'D2' '248.21'
Set-TextValue $ws 'D4' '5.536'
Set-TextValue $ws 'D5' '0.05620'
Set-TextValue $ws 'D6' '6.482'
Set-TextValue $ws 'D7' '0.8075'
Set-TextValue $ws 'D8' '1.047'
Set-TextValue $ws 'D9' '0.1433'
Set-TextValue $ws 'D10' '0.07331'
Set-TextValue $ws 'D11' '0.03209'
Set-TextValue $ws 'D12' '0.02919'
Set-TextValue $ws 'D13' '0.09263'
Set-TextValue $ws 'D14' '0.001672'
Set-TextValue $ws 'D15' '3.208'
Set-TextValue $ws 'D16' '0.04729'
Set-TextValue $ws 'D17' '0.0005813'
Set-TextValue $ws 'E17' '16OneONE'
Set-TextValue $ws 'D18' '0.006453'
Set-TextValue $ws 'D19' '0.005076'
Set-TextValue $ws 'D20' '0.001055'
Set-TextValue $ws 'D21' '0.0001500'
Set-TextValue $ws 'D22' '3.985'
Set-TextValue $ws 'D23' '3.380'
Set-TextValue $ws 'D25' '0.3318'
Set-TextValue $ws 'D26' '0.1255'
Set-TextValue $ws 'D27' '0.0003302'
Set-TextValue $ws 'D40' '0.04150'
Set-TextValue $ws 'D41' '0.003237'
Set-TextValue $ws 'E41' '40KickTokenKICKWorstin24h'
Set-TextValue $ws 'B42' 'BKEXToken'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D42' '0.1043'
Set-TextValue $ws 'E42' '41BKEXTokenBKK'
Set-TextValue $ws 'B43' 'CEJI'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws 'D43' '0.002971'
Set-TextValue $ws 'E43' '42CEJICEJI'
Set-TextValue $ws 'D44' '0.008585'
Set-TextValue $ws 'D45' '0.00005643'
Set-TextValue $ws 'D47' '0.6803'
Set-TextValue $ws 'D48' '0.01608'
